# Add a new slide ("Источники") before the final "Спасибо за внимание!" slide.
$p = $ppt.ActivePresentation

# The presentation currently has 7 slides; the last one (index 7) is
# "Спасибо за внимание!". Insert the new slide at position 7 (pushing the
# "Спасибо" slide to position 8) using the same "Title and Object" layout
# (CustomLayout #2 / slideLayout2.xml) that every other content slide uses.
$newSlide = $p.Slides.Add(7, 2)

# Title placeholder.
$titleShape = $newSlide.Shapes.Item(1)
$titleShape.Name = "Заголовок 1"
$titleShape.TextFrame.TextRange.Text = "Источники"

# Body / content placeholder with three hyperlinked source references.
$bodyShape = $newSlide.Shapes.Item(2)
$bodyShape.Name = "Объект 2"
$tr = $bodyShape.TextFrame.TextRange

$p1a = "https://"
$p1b = "www.kv.by/archive/index2007111801.htm"
$p1c = " - Аналоговый компьютер"
$url1 = "https://www.kv.by/archive/index2007111801.htm"

$p2a = "https://"
$p2b = "uic.vsu.ru/ccmuseum/comp/analog/index.htm"
$p2c = " - Аналоговые ЭВМ"
$url2 = "https://uic.vsu.ru/ccmuseum/comp/analog/index.htm"

$p3a = "https://analitikalmir.ru.gg/%"
$p3b = "26%231040%3B%26%231042%3B%26%231052%3B.htm"
$p3c = " – Аналоговая АВМ"
$url3 = "https://analitikalmir.ru.gg/%26%231040%3B%26%231042%3B%26%231052%3B.htm"

$tr.Text = $p1a + $p1b + $p1c + "`r" + $p2a + $p2b + $p2c + "`r" + $p3a + $p3b + $p3c

$pos = 1
$r1a = $tr.Characters($pos, $p1a.Length); $pos += $p1a.Length
$r1b = $tr.Characters($pos, $p1b.Length); $pos += $p1b.Length
$r1c = $tr.Characters($pos, $p1c.Length); $pos += $p1c.Length
$pos += 1

$r2a = $tr.Characters($pos, $p2a.Length); $pos += $p2a.Length
$r2b = $tr.Characters($pos, $p2b.Length); $pos += $p2b.Length
$r2c = $tr.Characters($pos, $p2c.Length); $pos += $p2c.Length
$pos += 1

$r3a = $tr.Characters($pos, $p3a.Length); $pos += $p3a.Length
$r3b = $tr.Characters($pos, $p3b.Length); $pos += $p3b.Length
$r3c = $tr.Characters($pos, $p3c.Length); $pos += $p3c.Length

$r1a.ActionSettings.Item(1).Hyperlink.Address = $url1
$r1b.ActionSettings.Item(1).Hyperlink.Address = $url1

$r2a.ActionSettings.Item(1).Hyperlink.Address = $url2
$r2b.ActionSettings.Item(1).Hyperlink.Address = $url2

$r3a.ActionSettings.Item(1).Hyperlink.Address = $url3
$r3b.ActionSettings.Item(1).Hyperlink.Address = $url3
